# Update social media content and assets
$d = $word.ActiveDocument

# 1) Intro paragraph text
$d.Content.Find.Execute(
    "Camp under the sky—well, museum‑style! Splash in the Arroyo stream, fish in blue‑sand lakes, hop trikes across a mini Golden Gate bridge, then circle up for campfire songs and crafts. It’s quintessential summer fun, perfectly scaled for little explorers and their families.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Spend an unforgettable night at Kidspace! Families will enjoy crafts, storytime, stargazing, and exclusive museum access. Bring your sleeping bags for a cozy camp experience under museum lights.",
    2
) | Out-Null

# 2) Date line
$d.Content.Find.Execute(
    "📅 Date: 2025‑07‑10–2025‑07‑13",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "📅 Date: 2025-07-18 – 2025-07-19",
    2
) | Out-Null

# 3) Time line
$d.Content.Find.Execute(
    "🕘 Time: Museum hours",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "🕘 Time: 6:00 pm – 8:00 am",
    2
) | Out-Null

# 4) Tickets line
$d.Content.Find.Execute(
    "💰 Tickets: Included",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "💰 Tickets: Varies",
    2
) | Out-Null

# 5) Hashtags line
$d.Content.Find.Execute(
    "#SummerCampout #NaturePlay #CampfireSongs #TrikeAdventure #SplashPlay #CreativeKids #FamilyFun #OutdoorSTEM #PasadenaEvents #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "#Campout #Kidspace #FamilyAdventure #MuseumSleepover #Stargazing #CraftNight #PasadenaKids #OvernightEvent #ScienceFun #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles",
    2
) | Out-Null
